$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark earlier Redux tasks as completed
$ws.Range("C8").Value = "Done"   # Products Back end
$ws.Range("C9").Value = "Done"   # Add Redux
$ws.Range("C10").Value = "Done"  # Add Redux To Products

# "Add Redux To Filter" is now the task being worked on
$ws.Range("C11").Value = "Active"

# Move selection to A6 as in the edited workbook
$ws.Range("A6").Select()
